$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D cells keep their original Text storage type (values are
# labels like "64.522.00" / "572.21", not numbers) by forcing the Text
# number format before assigning the new value.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.522.00'
$ws.Range('E2').Value = '  -3.03%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.173.12'
$ws.Range('E3').Value = '  -4.87%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '572.21'
$ws.Range('E5').Value = '  -2.63%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '169.68'
$ws.Range('E6').Value = '  -7.14%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.607'
$ws.Range('E7').Value = '  -6.60%  '
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.179.32'
$ws.Range('E9').Value = '  -4.52%  '
$ws.Range('E10').Value = '  -5.20%  '
$ws.Range('E11').Value = '  -0.24%  '
$ws.Range('E12').Value = '  -3.39%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.726.91'
$ws.Range('E13').Value = '  -4.69%  '
$ws.Range('E14').Value = '  -2.10%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '64.556.59'
$ws.Range('E15').Value = '  -2.88%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '25.45'
$ws.Range('E16').Value = '  -4.35%  '
$ws.Range('E17').Value = '  -4.60%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.180.38'
$ws.Range('E18').Value = '  -5.58%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '417.30'
$ws.Range('E19').Value = '  -2.20%  '
$ws.Range('E20').Value = '  -1.57%  '
$ws.Range('E21').Value = '  -3.74%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.15'
$ws.Range('E22').Value = '  -3.48%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  -0.08%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '70.19'
$ws.Range('E25').Value = '  -2.47%  '
$ws.Range('E26').Value = '  -0.24%  '
$ws.Range('E27').Value = '  -3.54%  '
$ws.Range('E28').Value = '  -8.73%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.81'
$ws.Range('E29').Value = '  -2.51%  '
$ws.Range('E30').Value = '  +0.11%  '
$ws.Range('E31').Value = '  -4.86%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '21.78'
$ws.Range('E32').Value = '  -3.05%  '
$ws.Range('E33').Value = '  -0.11%  '
$ws.Range('E34').Value = '  -3.02%  '
$ws.Range('E35').Value = '  -4.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.14'
$ws.Range('E36').Value = '  -4.99%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '155.94'
$ws.Range('E37').Value = '  -3.05%  '
$ws.Range('E38').Value = '  -5.39%  '
$ws.Range('E39').Value = '  -5.62%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.714.32'
$ws.Range('E40').Value = '  -5.52%  '
$ws.Range('E41').Value = '  -2.56%  '
$ws.Range('E42').Value = '  -8.23%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '39.09'
$ws.Range('E43').Value = '  -1.97%  '
$ws.Range('E44').Value = '  -5.86%  '
$ws.Range('E45').Value = '  -6.19%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.53'
$ws.Range('E46').Value = '  -7.50%  '
$ws.Range('E47').Value = '  -3.50%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '21.59'
$ws.Range('E48').Value = '  -7.30%  '
$ws.Range('B49').Value = 'Bittensor'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '292.97'
$ws.Range('E49').Value = '  -6.81%  '
$ws.Range('E50').Value = '  -13.24%  '
$ws.Range('E51').Value = '  -6.05%  '
